# Applies the "Finished animations for boss and integrated them into unity"
# edit to the Raid logic / Fire demon boss section of the document.
#
# Summary of the content change (see commit diff):
#   - "Raid logic" attack list loses its "Bone breaker" entry; the
#     remaining entries shift up and the first one is renamed.
#   - "Fire demon mesh" / "Animations" / the nested "Attack animations"
#     list and the trailing "Death" entry become green (color 92D050)
#     to match the rest of the "finished" / in-progress items.
#   - The nested attack-animation list also loses its "Bone breaker"
#     entry and gets its first entry renamed, mirroring the outer list.
#   - A new "Stun" animation entry is added after "Death".

$d = $word.ActiveDocument

$green = 5296274   # RGB(0x92, 0xD0, 0x50) == "92D050"

# --- "Raid logic" attack list (ilvl 3) ------------------------------------
# Normal attack, Bone breaker, 360 degree attack, Time based attack
# -> TurretDestroyer, 360 degree attack, Time based attack
$pNormalAttack = $d.Paragraphs.Item(22)
$pNormalAttack.Range.Text = "TurretDestroyer"

$pBoneBreaker = $d.Paragraphs.Item(23)
$pBoneBreaker.Range.Delete()

# --- Color "Fire demon mesh" through "Death" (now items 25-33) -----------
# Applied paragraph-by-paragraph so the paragraph-mark run properties
# (w:pPr/w:rPr) get the color too, matching how Word colors a selected
# block of whole paragraphs.
for ($i = 25; $i -le 33; $i++) {
    $pColor = $d.Paragraphs.Item($i)
    $pColor.Range.Font.Color = $green
}

# --- Nested "Attack animations" list (ilvl 4) -----------------------------
# Normal, Bone breaker, 360 degree attack, Time based attack
# -> TurretDestroyer, 360 degree attack, Time based attack
$pNormal = $d.Paragraphs.Item(29)
$pNormal.Range.Text = "TurretDestroyer"
$pNormal.Range.Font.Color = $green

$pBoneBreaker2 = $d.Paragraphs.Item(30)
$pBoneBreaker2.Range.Delete()

# --- Add "Stun" after "Death" (now item 32) -------------------------------
$pDeath = $d.Paragraphs.Item(32)
$pDeath.Range.InsertParagraphAfter()
$pStun = $d.Paragraphs.Item(33)
$pStun.Range.Text = "Stun"
$pStun.Range.Font.Color = $green
